# Weekly data refresh: a new week's record is inserted at the top of the
# "Camote" (Vega Modelo de Temuco) price table, pushing the existing rows
# down by one. The table previously ran from row 131 to row 189; after the
# edit it runs from row 131 to row 190 (dimension A1:R189 -> A1:R190).
#
# Concretely: insert a new row at 131 (shifting 131..189 down to 132..190,
# which is exactly what the diff shows - each old row's data reappears one
# row lower) and populate the newly-inserted row 131 with the latest week's
# record: same commodity / volume / price breakdown as the prior top row,
# but a new reporting date (serial 45027 = 2023-04-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data block (rows 131-189) down by one row.
$ws.Rows(131).Insert()

# Fill in the newly inserted row with this week's record.
$ws.Range("A131").Value = 10
$ws.Range("B131").Value = "Vega Modelo de Temuco"
$ws.Range("C131").Value = "La Araucanía"
$ws.Range("D131").Value = 45027
$ws.Range("E131").Value = 9
$ws.Range("F131").Value = 100114002
$ws.Range("G131").Value = "Camote"
$ws.Range("H131").Value = "Sin especificar"
$ws.Range("I131").Value = "Primera"
$ws.Range("J131").Value = 30
$ws.Range("K131").Value = 26000
$ws.Range("L131").Value = 26000
$ws.Range("M131").Value = 26000
$ws.Range("N131").Value = "$/malla 20 kilos"
$ws.Range("O131").Value = "Perú"
$ws.Range("P131").Value = 1300
$ws.Range("Q131").Value = 20
$ws.Range("R131").Value = "Hortaliza"
